# Auto-generated Excel COM-interop script
# Applies cached-market-price / leve-profit value updates across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 677.6
$ws.Range("J17").Value = 677.6
$ws.Range("L17").Value = 2032.8
$ws.Range("N17").Value = -2368.8
$ws.Range("H19").Value = 1270.5883
$ws.Range("I19").Value = 1258.8
$ws.Range("J19").Value = 1287.4286
$ws.Range("K19").Value = 1258.8
$ws.Range("L19").Value = 1287.4286
$ws.Range("M19").Value = -1083.8
$ws.Range("N19").Value = -1637.4286
$ws.Range("H32").Value = 3798.6667
$ws.Range("I32").Value = 3795
$ws.Range("J32").Value = 3800.5
$ws.Range("K32").Value = 3795
$ws.Range("L32").Value = 3800.5
$ws.Range("M32").Value = -3469
$ws.Range("N32").Value = -4452.5
$ws.Range("H98").Value = 1200.2858
$ws.Range("I98").Value = 880.6
$ws.Range("K98").Value = 880.6
$ws.Range("M98").Value = 617.4
$ws.Range("H106").Value = 5990
$ws.Range("I106").Value = 5980
$ws.Range("K106").Value = 5980
$ws.Range("M106").Value = -5349
$ws.Range("H107").Value = 1276.4286
$ws.Range("I107").Value = 1276.4286
$ws.Range("K107").Value = 1276.4286
$ws.Range("M107").Value = 643.5714
$ws.Range("H116").Value = 4208.4614
$ws.Range("I116").Value = 3861.75
$ws.Range("K116").Value = 3861.75
$ws.Range("M116").Value = -419.75
$ws.Range("H122").Value = 1200.2858
$ws.Range("I122").Value = 880.6
$ws.Range("K122").Value = 2641.8
$ws.Range("M122").Value = -191.8000000000002
$ws.Range("H129").Value = 1994
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 1994
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 5982
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -15982
$ws.Range("H132").Value = 2099.8
$ws.Range("J132").Value = 755.5
$ws.Range("L132").Value = 2266.5
$ws.Range("N132").Value = -7326.5
$ws.Range("H137").Value = 4988.6113
$ws.Range("I137").Value = 3123
$ws.Range("J137").Value = 5521.643
$ws.Range("K137").Value = 9369
$ws.Range("L137").Value = 16564.929
$ws.Range("M137").Value = -6819
$ws.Range("N137").Value = -21664.929
$ws.Range("H138").Value = 4875.095
$ws.Range("J138").Value = 5152.0527
$ws.Range("L138").Value = 15456.1581
$ws.Range("N138").Value = -25736.1581

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 2015
$ws.Range("I21").Value = 2015
$ws.Range("K21").Value = 2015
$ws.Range("M21").Value = -1641
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
$ws.Range("H32").Value = 2018.7142
$ws.Range("I32").Value = 1540.8298
$ws.Range("K32").Value = 1540.8298
$ws.Range("M32").Value = -1253.8298
$ws.Range("H45").Value = 1183.4
$ws.Range("I45").Value = 1095.5555
$ws.Range("J45").Value = 1974
$ws.Range("K45").Value = 1095.5555
$ws.Range("L45").Value = 1974
$ws.Range("M45").Value = -718.5554999999999
$ws.Range("N45").Value = -2728
$ws.Range("H61").Value = 2657.8823
$ws.Range("I61").Value = 2636.875
$ws.Range("K61").Value = 2636.875
$ws.Range("M61").Value = -2424.875
$ws.Range("H74").Value = 863.6667
$ws.Range("I74").Value = 863.6667
$ws.Range("K74").Value = 863.6667
$ws.Range("M74").Value = 10.33330000000001
$ws.Range("H77").Value = 863.6667
$ws.Range("I77").Value = 863.6667
$ws.Range("K77").Value = 4318.3335
$ws.Range("M77").Value = 49.66650000000027
$ws.Range("H136").Value = 2657.8823
$ws.Range("I136").Value = 2636.875
$ws.Range("K136").Value = 7910.625
$ws.Range("M136").Value = -5360.625

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 17644.777
$ws.Range("I20").Value = 17255.143
$ws.Range("J20").Value = 19008.5
$ws.Range("K20").Value = 17255.143
$ws.Range("L20").Value = 19008.5
$ws.Range("M20").Value = -17008.143
$ws.Range("N20").Value = -19502.5
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 69996
$ws.Range("J68").Value = 69996
$ws.Range("L68").Value = 69996
$ws.Range("N68").Value = -71494
$ws.Range("H71").Value = 69996
$ws.Range("J71").Value = 69996
$ws.Range("L71").Value = 209988
$ws.Range("N71").Value = -217476
$ws.Range("H86").Value = 48469.832
$ws.Range("I86").Value = 9477.75
$ws.Range("J86").Value = 126454
$ws.Range("K86").Value = 9477.75
$ws.Range("L86").Value = 126454
$ws.Range("M86").Value = -8354.75
$ws.Range("N86").Value = -128700
$ws.Range("H89").Value = 48469.832
$ws.Range("I89").Value = 9477.75
$ws.Range("J89").Value = 126454
$ws.Range("K89").Value = 47388.75
$ws.Range("L89").Value = 632270
$ws.Range("M89").Value = -41772.75
$ws.Range("N89").Value = -643502
$ws.Range("H94").Value = 1426.5
$ws.Range("J94").Value = 1370
$ws.Range("L94").Value = 1370
$ws.Range("N94").Value = -2272
$ws.Range("H122").Value = 1010
$ws.Range("I122").Value = 1010
$ws.Range("K122").Value = 3030
$ws.Range("M122").Value = -580
$ws.Range("H134").Value = 2423.2307
$ws.Range("I134").Value = 2423.2307
$ws.Range("K134").Value = 7269.6921
$ws.Range("M134").Value = -4734.6921

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 100
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H123").Value = 5828.6313
$ws.Range("I123").Value = 2744
$ws.Range("K123").Value = 8232
$ws.Range("M123").Value = -5782
$ws.Range("H132").Value = 3811.625
$ws.Range("I132").Value = 2000
$ws.Range("K132").Value = 18000
$ws.Range("M132").Value = -15470

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 40000
$ws.Range("J57").Value = 40000
$ws.Range("L57").Value = 40000
$ws.Range("N57").Value = -41640
$ws.Range("H80").Value = 6083.5713
$ws.Range("I80").Value = 5194
$ws.Range("J80").Value = 6750.75
$ws.Range("K80").Value = 5194
$ws.Range("L80").Value = 6750.75
$ws.Range("M80").Value = -4196
$ws.Range("N80").Value = -8746.75
$ws.Range("H83").Value = 6083.5713
$ws.Range("I83").Value = 5194
$ws.Range("J83").Value = 6750.75
$ws.Range("K83").Value = 25970
$ws.Range("L83").Value = 33753.75
$ws.Range("M83").Value = -20978
$ws.Range("N83").Value = -43737.75
$ws.Range("H126").Value = 1814

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1919.0714
$ws.Range("I22").Value = 1810.6364
$ws.Range("J22").Value = 2316.6667
$ws.Range("K22").Value = 1810.6364
$ws.Range("L22").Value = 2316.6667
$ws.Range("M22").Value = -1515.6364
$ws.Range("N22").Value = -2906.6667
$ws.Range("H27").Value = 1919.0714
$ws.Range("I27").Value = 1810.6364
$ws.Range("J27").Value = 2316.6667
$ws.Range("K27").Value = 1810.6364
$ws.Range("L27").Value = 2316.6667
$ws.Range("M27").Value = -1703.6364
$ws.Range("N27").Value = -2530.6667
$ws.Range("H46").Value = 3721.6667
$ws.Range("I46").Value = 3695
$ws.Range("K46").Value = 3695
$ws.Range("M46").Value = -3507
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H136").Value = 11640.1
$ws.Range("I136").Value = 12800.875
$ws.Range("K136").Value = 38402.625
$ws.Range("M136").Value = -35852.625

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 22726.773
$ws.Range("I13").Value = 22726.773
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 22726.773
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -22586.773
$ws.Range("N13").ClearContents()
$ws.Range("H31").Value = 5000
$ws.Range("J31").Value = 5000
$ws.Range("L31").Value = 5000
$ws.Range("N31").Value = -5696
$ws.Range("H136").Value = 8036.143
$ws.Range("I136").Value = 7566.2104
$ws.Range("K136").Value = 22698.6312
$ws.Range("M136").Value = -20148.6312
